$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D stores prices as plain-number-looking text (e.g. "242.94"). A bare
# Range.Value assignment of such a string auto-converts the cell to a numeric
# type in Excel, which would not match the source data (stored as text). So for
# every D-column cell we touch, force Text number format first, write the new
# price string, and afterwards reset the cell style back to "Normal" so the
# temporary Text format does not linger as a spurious formatting change.

# Row 2: update D
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.94"

# Row 3: update D
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.00"

# Row 4: update D
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.389"

# Row 5: update D
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05899"

# Row 6: update D
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.455"

# Row 7: update D
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.587"

# Row 9: update D
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9164"

# Row 10: update B, C, D, E
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01134"
$ws.Range("E10").Value = "9OneONEBestin24h"

# Row 11: update B, C, D, E
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1424"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12: update B, C, D, E
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07392"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13: update B, C, D, E
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03254"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14: update B, C, D, E
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03059"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15: update B, C, D, E
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09339"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16: update B, C, D, E
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.861"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17: update B, C, D, E
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001565"
$ws.Range("E17").Value = "16BitForexTokenBF"

# Row 18: update B, C, D, E
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04662"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19: update D
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005906"

# Row 20: update D, E
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001285"
$ws.Range("E20").Value = "19BitKanKAN"

# Row 21: update D
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004898"

# Row 22: update D
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009504"

# Row 23: update D
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.618"

# Row 24: update D
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.150"

# Row 25: update D
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3231"

# Row 41: update D
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006182"

# Row 42: update D
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1070"

# Row 43: update D
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003001"

# Row 44: update D
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008092"

# Row 45: update D
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005187"

# Row 47: update D
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7503"

# Row 48: update D
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002279"

# Restore default styling on every price cell we reformatted above.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
